# Applies the "Birthday and Ouput Fixed" edits to the 何泽恩 enrollment form.
# Several old values are duplicated verbatim elsewhere in the document
# (e.g. phone numbers, ID numbers, "否", "无", "111"), so instead of a
# blind document-wide Find/Replace we locate the Nth occurrence of each
# search string (counted from the top of the document) and replace only
# that Range's text. Occurrence indices below were determined by
# inspecting the document's table/cell layout so each edit lands on the
# exact field the diff touches. NOTE: this runtime's PowerShell only
# binds positional parameters, so Replace-NthMatch is always called
# positionally (Needle, Occurrence, NewText).

$d = $word.ActiveDocument

function Replace-NthMatch {
    param([string]$Needle, [int]$Occurrence, [string]$NewText)

    $rng = $d.Content
    $rng.Start = 0
    $rng.End = $d.Content.End

    $count = 0
    $target = $null
    while ($rng.Find.Execute($Needle)) {
        $count++
        if ($count -eq $Occurrence) {
            $target = $d.Range($rng.Start, $rng.End)
            break
        }
        $rng.Start = $rng.End
        $rng.End = $d.Content.End
    }

    if ($target -eq $null) {
        throw "Replace-NthMatch: could not find occurrence $Occurrence of $Needle"
    }

    $target.Text = $NewText
}

# 1. Birth date: 2020年10月28日 -> 2020年10月29日 (unique)
Replace-NthMatch "2020年10月28日" 1 "2020年10月29日"

# 2. Ethnicity: 蒙古族 -> 维吾尔族 (unique)
Replace-NthMatch "蒙古族" 1 "维吾尔族"

# 3. Father's work unit: 11111 -> 父亲工作单位 (unique full match; do before the
#    "111" replacement below since "111" is a substring of "11111")
Replace-NthMatch "11111" 1 "父亲工作单位"

# 4. Home address: 111 -> 家庭住址 (now unique after step 3 removed the other
#    "111"-containing value)
Replace-NthMatch "111" 1 "家庭住址"

# 5. Grandparent's name (4th occurrence of 何泽恩 in the document) -> empty
Replace-NthMatch "何泽恩" 4 ""

# 6. Grandparent's work unit: 祖辈工作单位 -> empty (unique)
Replace-NthMatch "祖辈工作单位" 1 ""

# 7. Grandparent's phone (3rd occurrence of 19851937930) -> empty
Replace-NthMatch "19851937930" 3 ""

# 8. "是否爱提问" answer (2nd occurrence of 否, the standalone value cell) -> 是
Replace-NthMatch "否" 2 "是"

# 9. "是否有午睡的习惯" answer -> 是. After step 8 replaced the old
#    occurrence 2, the remaining value cell (previously occurrence 5) is
#    now occurrence 4.
Replace-NthMatch "否" 4 "是"

# 10. Interests/hobbies tags -> 篮球,足球,跑步,羽毛球 (unique)
Replace-NthMatch "样例标签1（请删除）,样例标签2（请删除）" 1 "篮球,足球,跑步,羽毛球"

# 11. Morning wake-up time: 06:10 -> 06:15 (unique)
Replace-NthMatch "06:10" 1 "06:15"

# 12. Evening bedtime: 20:00 -> 22:00 (unique)
Replace-NthMatch "20:00" 1 "22:00"

# 13. "有何病史" answer (1st occurrence of 无) -> 双方都是
Replace-NthMatch "无" 1 "双方都是"

# 14. Drug allergy: 头孢 -> 111 (unique)
Replace-NthMatch "头孢" 1 "111"

# 15. Special-needs note: 特殊情况 -> empty (unique)
Replace-NthMatch "特殊情况" 1 ""

# 16. Father's ID number (3rd occurrence of 460035200101121517) -> 46003519760820152X
Replace-NthMatch "460035200101121517" 3 "46003519760820152X"
